$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class "0")
$ws.Range("B2").Value = 0.8260869565217391
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.9047619047619047

# Row 3 (class "1")
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.8545454545454545
$ws.Range("D3").Value = 0.9215686274509803

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9139784946236559
$ws.Range("C4").Value = 0.9139784946236559
$ws.Range("D4").Value = 0.9139784946236559
$ws.Range("E4").Value = 0.9139784946236559

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9130434782608696
$ws.Range("C5").Value = 0.9272727272727272
$ws.Range("D5").Value = 0.9131652661064424

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9289387564282375
$ws.Range("C6").Value = 0.9139784946236559
$ws.Range("D6").Value = 0.9147013644167344
